$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numbers of Commits in GitHub increased for both team members -
# the dependent "Total Score" cell (C44, a SUM formula) recalculates
# automatically.
$ws.Range("C8").Value = 13
$ws.Range("C9").Value = 16
